# Append the 2025-11-09 profit allocation row to the sheet (row 69).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to be treated as plain text so the date-like string
# "11/09/2025" is stored as text (matching the other Date column cells)
# instead of being auto-converted into a date serial number.
$ws.Range("A69").NumberFormat = "@"
$ws.Range("A69").Value = "11/09/2025"
# Reset the style back to the default/Normal style so no extra
# number-format styling is left on the new cell (keeps it consistent
# with the rest of the Date column, which uses the default style).
$ws.Range("A69").Style = "Normal"

$ws.Range("B69").Value = 0.1958362470302913
$ws.Range("C69").Value = 0.8041637529697087
